# working mobile rev 1.0
#
# 1) On sheet "10jul2022", update B19 from 0 -> 5.
# 2) Add a new sheet "10jun2022" (after "10jul2022") populated with the
#    Point/Riesgo/NDVI/Latitud/Longitud table (rows 1-25), reusing the
#    same header style as the existing sheets.

$wb = $excel.ActiveWorkbook

# --- 1) Fix the single value on the existing "10jul2022" sheet ---
$wsJul = $wb.Worksheets.Item("10jul2022")
$wsJul.Range("B19").Value = 5

# --- 2) Build the new "10jun2022" sheet at the end of the workbook ---
$wsFirst = $wb.Worksheets.Item(1)
$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsLast)
$wsNew.Name = "10jun2022"

# Header row
$wsNew.Range("A1").Value = "Point"
$wsNew.Range("B1").Value = "Riesgo"
$wsNew.Range("C1").Value = "NDVI"
$wsNew.Range("D1").Value = "Latitud"
$wsNew.Range("E1").Value = "Longitud"

# Match the bold/centered/bordered header style used on the other sheets
$wsFirst.Range("A1:E1").Copy()
$wsNew.Range("A1:E1").PasteSpecial(-4122)

# Data rows (Point, Riesgo, NDVI, Latitud, Longitud)
$data = @(
    @(1, 0, 0.8354177186702311, 6.15818348129518, -75.78945109963882),
    @(2, 0, 0.8054641585934869, 6.158839773883416, -75.78945109963882),
    @(3, 0, 0.7842966904325106, 6.159496066471652, -75.78945109963882),
    @(4, 0, 0.8078793282470147, 6.160152359059888, -75.78945109963882),
    @(5, 0, 0.7468479526476191, 6.160808651648123, -75.78945109963882),
    @(6, 0, 0.6688650456038993, 6.15818348129518, -75.7885922824137),
    @(7, 0, 0.803323712101627, 6.158839773883416, -75.7885922824137),
    @(8, 0, 0.7724477892480791, 6.159496066471652, -75.7885922824137),
    @(9, 0, 0.8173297342579228, 6.160152359059888, -75.7885922824137),
    @(10, 0, 0.8754531743872137, 6.160808651648123, -75.7885922824137),
    @(11, 0, 0.6657485699358098, 6.15818348129518, -75.78773346518857),
    @(12, 0, 0.8600671006557239, 6.158839773883416, -75.78773346518857),
    @(13, 4, 0.8757465203616336, 6.159496066471652, -75.78773346518857),
    @(14, 0, 0.8046835140365085, 6.160152359059888, -75.78773346518857),
    @(15, 0, 0.8830674665163043, 6.160808651648123, -75.78773346518857),
    @(16, 0, 0.4625837282354751, 6.15818348129518, -75.78687464796343),
    @(17, 0, 0.7907732269250011, 6.158839773883416, -75.78687464796343),
    @(18, 0, 0.7704232451339251, 6.159496066471652, -75.78687464796343),
    @(19, 0, 0.804655329489094, 6.160152359059888, -75.78687464796343),
    @(20, 0, 0.875232071071354, 6.160808651648123, -75.78687464796343),
    @(21, 0, 0.4250888714528885, 6.15818348129518, -75.78601583073831),
    @(22, 0, 0.5670474935648414, 6.158839773883416, -75.78601583073831),
    @(23, 0, 0.2676151827328355, 6.159496066471652, -75.78601583073831),
    @(24, 0, 0.7720649349443865, 6.160152359059888, -75.78601583073831),
    @(25, 0, 0.8234034213899299, 6.160808651648123, -75.78601583073831)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $wsNew.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex = $colIndex + 1
    }
    $rowIndex = $rowIndex + 1
}

# Restore the original active sheet / selection so the workbook-level view
# state matches the source file (first sheet active).
$wsNew.Range("A1").Select()
$wsFirst.Activate()
